$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.222.90"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.596.91"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.Value = "'212.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$c = $ws.Range("D6")
$c.Value = "'0.503"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  +0.41%  "
$c = $ws.Range("D10")
$c.Value = "'18.93"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "
$c = $ws.Range("D11")
$c.Value = "'0.0853"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "1.820.94"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.602.42"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("E15").Value = "  -2.22%  "
$c = $ws.Range("D16")
$c.Value = "'63.63"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "26.206.06"
$ws.Range("E17").Value = "  +0.05%  "
$c = $ws.Range("D18")
$c.Value = "'228.09"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +6.47%  "
$ws.Range("E19").Value = "  -0.74%  "
$c = $ws.Range("D20")
$c.Value = "'7.56"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.10%  "
$c = $ws.Range("D23")
$c.Value = "'2.17"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$c = $ws.Range("D24")
$c.Value = "'8.89"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.66%  "
$c = $ws.Range("D25")
$c.Value = "'145.47"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E29").Value = "  +1.51%  "
$c = $ws.Range("D30")
$c.Value = "'0.0493"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "1.446.63"
$ws.Range("E33").Value = "  +4.14%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("E38").Value = "  -1.54%  "
$c = $ws.Range("D40")
$c.Value = "'5.74"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  +1.77%  "
$c = $ws.Range("D43")
$c.Value = "'0.929"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").Value = "1.734.10"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -1.11%  "
$c = $ws.Range("D46")
$c.Value = "'60.35"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.99%  "
$c = $ws.Range("D47")
$c.Value = "'87.50"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.78%  "
$c = $ws.Range("D48")
$c.Value = "'1.48"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D51")
$c.Value = "'7.42"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
